$wb = $excel.ActiveWorkbook

# Both the "展览" sheet and the "全部类型" sheet carry the same table of
# convention data and both need the same refreshed numbers.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # row 2: 想去人数 25 -> 29, 最低票价 "59" -> "已停售"
    $ws.Range("F2").Value = 29
    $ws.Range("G2").Value = "已停售"

    # row 3: 想去人数 1835 -> 1858
    $ws.Range("F3").Value = 1858

    # row 4: 想去人数 557 -> 566
    $ws.Range("F4").Value = 566

    # row 5: 想去人数 1185 -> 1219
    $ws.Range("F5").Value = 1219

    # row 6: 想去人数 6142 -> 6232
    $ws.Range("F6").Value = 6232

    # row 7: 想去人数 146 -> 151
    $ws.Range("F7").Value = 151
}
